$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.449.15'
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('D3').Value = '1.935.26'
$ws.Range('E3').Value = '  -1.68%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '242.01'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.20%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.608'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.89%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '56.27'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -5.28%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.358'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -4.66%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0809'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -4.47%  '
$ws.Range('E11').Value = '  -1.34%  '
$ws.Range('D12').Value = '2.220.82'
$ws.Range('E12').Value = '  -1.59%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '21.03'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -5.40%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.802'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -4.71%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '13.22'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.84%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.12'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -4.74%  '
$ws.Range('D17').Value = '1.934.00'
$ws.Range('E17').Value = '  -2.46%  '
$ws.Range('D18').Value = '36.346.81'
$ws.Range('E18').Value = '  -0.44%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '68.90'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.40%  '
$ws.Range('D20').Value = '0.0₃0851'
$ws.Range('E20').Value = '  -3.67%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '226.35'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.39%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.93'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.77%  '
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.38'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -5.84%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.28'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.05'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -5.71%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '159.16'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -3.89%  '
$ws.Range('E28').Value = '  +7.74%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.03'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -3.80%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.117'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.37%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.08'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -8.08%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.54'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -5.32%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0609'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -5.13%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.10'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -6.00%  '
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.08'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.36%  '
$ws.Range('E37').Value = '  -1.02%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.16'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.18%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.11'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +6.23%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0982'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.18%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.91'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.30%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0207'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.03%  '
$ws.Range('E43').Value = '  -5.33%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '15.59'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.58%  '
$ws.Range('D45').Value = '1.328.36'
$ws.Range('E45').Value = '  -1.80%  '
$ws.Range('E46').Value = '  -3.83%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '85.41'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.60%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.09'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -4.29%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.80'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.92%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.52'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +14.50%  '
$ws.Range('D51').Value = '2.113.60'
$ws.Range('E51').Value = '  -1.58%  '
